$d = $word.ActiveDocument
$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($p.Range.Text -like "*Modellverständlichkeit und Akzeptanz*") {
        $r = $p.Range
        $r.Collapse(0)
        $r.InsertParagraphAfter()
        break
    }
}
$newIdx = $idx + 1
$newPara = $d.Paragraphs($newIdx)
$nr = $newPara.Range
$titleStart = $nr.Start
$title = "Fehlendes Data Dictionary"
$nr.InsertAfter($title)
$titleEnd = $titleStart + $title.Length

# Use the paragraph's own range (live) to bold the title - it currently == title only
$titleRange = $newPara.Range
Write-Output "titleRange text=[$($titleRange.Text)]"
$titleRange.Font.Bold = $true
$titleRange.Font.BoldBi = $true

$nr3 = $newPara.Range
$breakPos = $titleEnd
$nr3.InsertBreak(6)

$bodyStart = $breakPos + 1
$body = "Da kein Data Dictionary zur Verfügung steht, besteht die Möglichkeit, dass bestimmte Variablen oder Werte nicht die erwartete Bedeutung haben oder falsch interpretiert werden. Insbesondere besteht das Risiko, dass zentrale Kennzahlen wie beispielsweise „damage“ betriebswirtschaftlich nicht sinnvoll definiert sind. Eine Klärung ist daher frühzeitig vorzunehmen."
$nr5 = $newPara.Range
$nr5.InsertAfter($body)
$bodyEnd = $bodyStart + $body.Length

$bodyRange = $d.Range($bodyStart, $bodyEnd)
Write-Output "bodyRange text=[$($bodyRange.Text)]"
$bodyRange.Bold = 0
$bodyRange.BoldBi = 0
Write-Output "bodyRange Bold=$($bodyRange.Bold) BoldBi=$($bodyRange.BoldBi)"
